$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login credentials used in rows 2-4 (A2:A4 / B2:B4)
# Set column B first so the new shared string "sEvYmEq" is inserted
# before "mngr601962", matching the expected shared string order.
$ws.Range("B2:B4").Value = "sEvYmEq"
$ws.Range("A2:A4").Value = "mngr601962"

# Update the selected range shown in the sheet view
$ws.Range("A2:B4").Select()
